$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.026.93"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "'2.540.91"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'544.10"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").Value = "'145.29"
$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("D9").Value = "'2.573.94"
$ws.Range("E9").Value = "  +2.28%  "

$ws.Range("D10").Value = "'0.102"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").Value = "'5.55"
$ws.Range("E12").Value = "  +2.65%  "

$ws.Range("E13").Value = "  +1.28%  "

$ws.Range("D14").Value = "'2.988.55"
$ws.Range("E14").Value = "  +2.02%  "

$ws.Range("D15").Value = "'23.92"
$ws.Range("E15").Value = "  -2.09%  "

$ws.Range("D16").Value = "'59.923.38"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("E17").Value = "  +2.57%  "

$ws.Range("D18").Value = "'2.549.50"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("D19").Value = "'11.37"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").Value = "'4.34"
$ws.Range("E20").Value = "  -0.89%  "

$ws.Range("D21").Value = "'328.50"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").Value = "'5.94"
$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("D24").Value = "'62.74"
$ws.Range("E24").Value = "  +2.05%  "

$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").Value = "'0.990"
$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("D28").Value = "'8.03"
$ws.Range("E28").Value = "  +2.17%  "

$ws.Range("D29").Value = "'7.13"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").Value = "'0.0₃0802"
$ws.Range("E30").Value = "  +1.45%  "

$ws.Range("E31").Value = "  -1.04%  "

$ws.Range("D32").Value = "'1.24"
$ws.Range("E32").Value = "  -3.85%  "

$ws.Range("D33").Value = "'163.12"
$ws.Range("E33").Value = "  +2.35%  "

$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  +4.85%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "'18.81"
$ws.Range("E36").Value = "  +0.37%  "

$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").Value = "'1.64"
$ws.Range("E38").Value = "  -2.14%  "

$ws.Range("D39").Value = "'5.77"
$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("D40").Value = "'37.24"
$ws.Range("E40").Value = "  +1.51%  "

$ws.Range("D41").Value = "'304.04"
$ws.Range("E41").Value = "  -3.70%  "

$ws.Range("D42").Value = "'0.843"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("D43").Value = "'3.75"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.610"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.993"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").Value = "'10.86"
$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("D47").Value = "'19.10"
$ws.Range("E47").Value = "  +2.34%  "

$ws.Range("D48").Value = "'0.0940"
$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").Value = "'125.08"
$ws.Range("E49").Value = "  -0.38%  "

$ws.Range("D50").Value = "'0.0522"
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("E51").Value = "  -0.98%  "
